$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Insert a new column before column B, shifting old B->C and old C->D.
$ws.Columns("B").Insert()

# New header text for the inserted column.
$ws.Range("B1").Value = "StatQuery"

# New query text for the inserted column.
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.gender IN ['FEMALE']  OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Match the wrap-text style used by A2 on the new B2 cell.
$ws.Range("B2").WrapText = $true

# Column B should be the same width as column A (75.81640625 characters).
$ws.Range("B1").ColumnWidth = 75

# Update selection to reflect the saved view state.
$ws.Range("A4").Select()
